# Update results & change main.py default argument
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B13: default argument in main.py changed -> run description now mentions "MLP 3 frames" ---
$ws.Range("B13").Value = "PPO use step distance reward + multiply critic lr + train every episode + MLP 3 frames vs. Random"

# --- New "Run dir" column (O) aligned with the existing rows ---
$ws.Range("O1").Value = "Run dir"
$ws.Range("O2").Value = "baseline_map*"
$ws.Range("O3").Value = "map*_use_dist"
$ws.Range("O4").Value = "map*_use_hit_wall"
$ws.Range("O5").Value = "map*_use_hit_wall_[actor]_2layers"
$ws.Range("O6").Value = "map*_use_hit_wall_[MLP]_[frames]3"
$ws.Range("O7").Value = "map*_use_hit_wall_[MLP]_[frames]9"
$ws.Range("O8").Value = "map*_use_hit_wall_[CNN]_[frames]3"
$ws.Range("O9").Value = "map*_use_hit_wall_[CNN]_[frames]9"
$ws.Range("O10").Value = "map*_use_summed_hit_wall_[MLP]_[frames]3"
$ws.Range("O11").Value = "map*_use_summed_hit_wall_[MLP]_[frames]9"
$ws.Range("O13").Value = "map*_use_step_dist"

# --- New row 12 (new observations), with fill colours matching existing "works"/"meh" notes ---
$ws.Range("J12").Value = "经常折返"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("J12").PasteSpecial(-4122) | Out-Null

$ws.Range("K12").Value = "不是稳赢，经常折返"
$ws.Range("I10").Copy() | Out-Null
$ws.Range("K12").PasteSpecial(-4122) | Out-Null

# --- New C13 note, same "works ok" fill as other green cells ---
$ws.Range("C13").Value = "反而是换边更顺，不过都还行"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Column widths manually widened for the new columns N/O ---
$ws.Columns.Item(14).ColumnWidth = 16.5
$ws.Columns.Item(15).ColumnWidth = 20.5

# --- View state: zoom + new selection ---
$excel.ActiveWindow.Zoom = 110
$ws.Range("I11").Select() | Out-Null
